# excel read and write file
# Update a few cells in the test fruit data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# price for Mango changed 299 -> 350
$ws.Range("D2").Value = 350

# fruit_name for row 3 changed Apple -> Iphone
$ws.Range("B3").Value = "Iphone"

# fruit_name for row 5 changed Banana -> Replace text
$ws.Range("B5").Value = "Replace text"
